$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8368080854415894
$ws.Range("B1").Value = 2.103044033050537
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.434626579284668
$ws.Range("E1").Value = 0.4855068325996399
